# Delete row 205 ("「安らぎを感じますか？」..." entry) from the active worksheet.
# This shifts all subsequent rows up by one, which matches the diff:
#   - dimension shrinks from A1:C292 to A1:C291
#   - every row previously numbered 206..292 becomes 205..291 with identical content
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(205).Delete()

$wb.Save()
